# Commit: "150 Less is more"
#
# The dashboard sheet lists several "blog" widget configs as cell text in
# row 11 (B11, D11, I11). I11 held the placeholder for blog post "ser: 148";
# it is being replaced with "ser: 150" (hence "150 Less is more" - one
# fewer placeholder, 150 now represented instead of 148).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I11").Value = "type: blog`nwidth: 2`nheight: 1`nser: 150"
